$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 fresh rows at 50-54 (pushing the old, empty trailing row down) so the
# new rows pick up the same "A/B/D/E styled" cell formatting that every other
# data row in the sheet already uses, then drop the now-empty row that got
# displaced to the bottom so the used range stays A1:G54.
$ws.Rows.Item(50).Insert()
$ws.Rows.Item(51).Insert()
$ws.Rows.Item(52).Insert()
$ws.Rows.Item(53).Insert()
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(55).Delete()

# --- New rows 50-53: completeme/COMPLETEME, amazy/AMAZY, tkanan/TKANAN, cmulia/CMULIA ---
# Column B (fullname, uppercase) is entered first for each row, then column A
# (username, lowercase), then the remaining columns (password / badan_usaha /
# divisi / role / approval) are filled in per row.
$ws.Cells.Item(50, 2).Value = "COMPLETEME"
$ws.Cells.Item(51, 2).Value = "AMAZY"
$ws.Cells.Item(52, 2).Value = "TKANAN"
$ws.Cells.Item(53, 2).Value = "CMULIA"

$ws.Cells.Item(50, 1).Value = "completeme"
$ws.Cells.Item(51, 1).Value = "amazy"
$ws.Cells.Item(52, 1).Value = "tkanan"
$ws.Cells.Item(53, 1).Value = "cmulia"

$ws.Cells.Item(50, 3).Value = "complete123"
$ws.Cells.Item(50, 4).Value = "COMPLETEME"
$ws.Cells.Item(50, 5).Value = "COMPLETEME"
$ws.Cells.Item(50, 6).Value = "USER"
$ws.Cells.Item(50, 7).Value = "ADMIN"

$ws.Cells.Item(51, 3).Value = "complete123"
$ws.Cells.Item(51, 4).Value = "AMAZY"
$ws.Cells.Item(51, 5).Value = "AMAZY"
$ws.Cells.Item(51, 6).Value = "USER"
$ws.Cells.Item(51, 7).Value = "ADMIN"

$ws.Cells.Item(52, 3).Value = "complete123"
$ws.Cells.Item(52, 4).Value = "TKANAN"
$ws.Cells.Item(52, 5).Value = "TKANAN"
$ws.Cells.Item(52, 6).Value = "USER"
$ws.Cells.Item(52, 7).Value = "ADMIN"

$ws.Cells.Item(53, 3).Value = "complete123"
$ws.Cells.Item(53, 4).Value = "CMULIA"
$ws.Cells.Item(53, 5).Value = "CMULIA"
$ws.Cells.Item(53, 6).Value = "USER"
$ws.Cells.Item(53, 7).Value = "ADMIN"

# --- New row 54: new divisi (AUDIT), area/username "mega", role EXECUTOR ---
$ws.Cells.Item(54, 1).Value = "mega"
$ws.Cells.Item(54, 5).Value = "AUDIT"
$ws.Cells.Item(54, 2).Value = "AUDITOR"
$ws.Cells.Item(54, 6).Value = "EXECUTOR"
$ws.Cells.Item(54, 3).Value = "complete123"
$ws.Cells.Item(54, 4).Value = "CV.CS"
$ws.Cells.Item(54, 7).Value = "ADMIN"

# --- View state: new active selection below the added rows ---
$ws.Range("D56").Select()
